# "Zeitblatt auf Stand gebracht" - bring the Jim Frey Zeitblatt workbook up to date.
#
# Summary of the real-world edit:
#  - Oktober sheet: log 2h on the first day (CCS reinstall / Tiva I2C code) and
#    2.5h on the second day (Tiva I2C libs, first build, datasheet reading),
#    with a short note in the Details column (wrapped for the 2nd entry).
#  - Jänner sheet: the "Stunden" (B) column got formatted for the whole month
#    (ready for hour entries) and the sheet got a portrait/A4 page setup.
#  - Selections / active sheet moved around as the author worked through it.

$wb = $excel.ActiveWorkbook

$oktober = $wb.Worksheets.Item("Oktober")
$jaenner = $wb.Worksheets.Item("Jänner")

# --- Oktober: fill in the first two days of actual work -------------------

$oktober.Cells.Item(4, 2).NumberFormat = "0.00"
$oktober.Cells.Item(4, 2).Value = 2
$oktober.Cells.Item(4, 3).Value = 'Neu Installieren von CCS, versucht Tiva Controller I2C Code zu gernerieren '

$oktober.Cells.Item(5, 2).NumberFormat = "0.00"
$oktober.Cells.Item(5, 2).Value = 2.5
$oktober.Cells.Item(5, 3).WrapText = $true
$oktober.Cells.Item(5, 3).Value = 'Tiva I2C Libs installiert, Projekt erstellt, Kompiliert mit includes aber beim Aufruf von Methoden der Includes gibt’s Fehler. Lichtesensor Datasheet duchgearbeitet - Paket funktion usw. '
$oktober.Rows.Item(5).RowHeight = 47.25

# --- Jänner: pre-format the Stunden column for the whole month, add a page setup

for ($r = 4; $r -le 34; $r++) {
    $jaenner.Cells.Item($r, 2).NumberFormat = "0.00"
}
$jaenner.Cells.Item(5, 3).WrapText = $true

$jaenner.PageSetup.PaperSize = 9
$jaenner.PageSetup.Orientation = 1

# --- selections: Jänner was being looked at, then work moved to Oktober ---

$jaenner.Activate() | Out-Null
$jaenner.Range("C9").Select() | Out-Null

$oktober.Activate() | Out-Null
$oktober.Range("B23").Select() | Out-Null
